$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts old rows 2-4 down to 3-5)
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the bold header row).
# Reset it to the plain/default style used by the other data rows, then re-apply
# the date style to column D to match the other date cells.
$ws.Rows.Item(2).Style = "Normal"

# Populate the new row 2 with data (mirrors columns A,B,C,E,F,G,H,I,N,Q,R from the row below it,
# and introduces new values for D,J,K,L,M,O,P)
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44643
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 100112052
$ws.Cells.Item(2, 7).Value = "Albahaca"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 90
$ws.Cells.Item(2, 11).Value = 2800
$ws.Cells.Item(2, 12).Value = 3000
$ws.Cells.Item(2, 13).Value = 2911
$ws.Cells.Item(2, 14).Value = "$/docena de matas"
$ws.Cells.Item(2, 15).Value = "Región Metropolitana"
$ws.Cells.Item(2, 16).Value = 485
$ws.Cells.Item(2, 17).Value = 6
$ws.Cells.Item(2, 18).Value = "Hortaliza"

# Match the date-number-format style used by the other rows' Fecha (column D) cells
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat
